# COREESG_holdings.xlsx - "Add files via upload" commit replay
#
# The workbook ships with a protected worksheet (legacy password hash
# "D382"). Excel's object model refuses any cell write while a sheet is
# protected, so we temporarily unprotect, apply the data refresh, then
# re-protect with the same password so the sheet is left in the same
# (protected) state it started in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$wasProtected = $ws.ProtectContents
if ($wasProtected) {
    $ws.Unprotect("D382")
}

# --- Update the "as of" date in the confidential disclosure banner ---
$ws.Range("A10").Value2 = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-26 for illustrative purposes only and are subject to change."

# --- Refresh the weight / percent-change figures (columns D & E) ---
$ws.Range("D2").Value = 0.2489687230967645
$ws.Range("E2").Value = 0.001809508142786687

$ws.Range("D3").Value = 0.4972942692058721
$ws.Range("E3").Value = 0.001843560705820257

$ws.Range("D4").Value = 0.09698469093595484
$ws.Range("E4").Value = 0.002149277061352262

$ws.Range("D5").Value = 0.1008555395951505
$ws.Range("E5").Value = 0.006294471811713276

$ws.Range("D6").Value = 0.05589677716625818
$ws.Range("E6").Value = 0.01531078610603287

$ws.Range("E7").Value = 0.003066406027457447

if ($wasProtected) {
    $ws.Protect("D382")
}
